$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid (G) and Absent (H) -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count (D) and Real (E) -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: D and E -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent (H) -> 1
$ws.Range("H6").Value = 1

# Row 7: H -> 1
$ws.Range("H7").Value = 1

# Row 8: H -> 1
$ws.Range("H8").Value = 1

# Row 9: D and E -> 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: H -> 1
$ws.Range("H10").Value = 1

# Row 11: H -> 1
$ws.Range("H11").Value = 1

# Row 12: H -> 1
$ws.Range("H12").Value = 1

# Row 13: D and E -> 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: H -> 1
$ws.Range("H14").Value = 1

# Row 15: H -> 1
$ws.Range("H15").Value = 1

# Row 16: H -> 1
$ws.Range("H16").Value = 1

# Row 17: D and E -> 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18: H -> 1
$ws.Range("H18").Value = 1
